$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last refreshed" timestamp label ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Mayo de 2020 a las 14:35"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 1667437
$ws.Range("C4").Value = 609
$ws.Range("E4").Value = 1121819
$ws.Range("G4").Value = 8
$ws.Range("H4").Value = 98691

# --- Row 11: Alemania ---
$ws.Range("B11").Value = 179992
$ws.Range("C11").Value = 6
$ws.Range("E11").Value = 11326

# --- Row 31: Portugal ---
$ws.Range("B31").Value = 30623
$ws.Range("C31").Value = 152
$ws.Range("D31").Value = 17549
$ws.Range("E31").Value = 11758
$ws.Range("G31").Value = 14
$ws.Range("H31").Value = 1316

# --- Rows 70/71: Luxemburgo & Azerbaiyan swap ranking (Azerbaiyan moves ahead) ---
$ws.Range("A70").Value = "Azerbaiyan"
$ws.Range("B70").Value = 4122
$ws.Range("C70").Value = 140
$ws.Range("D70").Value = 2607
$ws.Range("E70").Value = 1466
$ws.Range("H70").Value = 49

$ws.Range("A71").Value = "Luxemburgo"
$ws.Range("B71").Value = 3990
$ws.Range("C71").Value = 0
$ws.Range("D71").Value = 3758
$ws.Range("E71").Value = 123
$ws.Range("H71").Value = 109

# --- Row 86: Croacia ---
$ws.Range("B86").Value = 2244
$ws.Range("C86").Value = 1
$ws.Range("D86").Value = 2027
$ws.Range("E86").Value = 118

# --- Rows 103/104: Guinea-Bisau & Sri Lanka swap ranking (Sri Lanka moves ahead) ---
$ws.Range("A103").Value = "Sri Lanka"
$ws.Range("B103").Value = 1117
$ws.Range("C103").Value = 28
$ws.Range("D103").Value = 674
$ws.Range("E103").Value = 434
$ws.Range("H103").Value = 9

$ws.Range("A104").Value = "Guinea-Bisau"
$ws.Range("B104").Value = 1114
$ws.Range("C104").Value = 0
$ws.Range("D104").Value = 42
$ws.Range("E104").Value = 1066
$ws.Range("H104").Value = 6

# --- Row 105: Libano ---
$ws.Range("B105").Value = 1114
$ws.Range("C105").Value = 17
$ws.Range("D105").Value = 688
$ws.Range("E105").Value = 400

# --- Row 134: Madagascar ---
$ws.Range("D134").Value = 142
$ws.Range("E134").Value = 383

# --- Rows 157/158: Islas Feroe & Mozambique swap ranking (Mozambique moves ahead) ---
$ws.Range("A157").Value = "Mozambique"
$ws.Range("B157").Value = 194
$ws.Range("C157").Value = 26
$ws.Range("D157").Value = 48
$ws.Range("E157").Value = 146

$ws.Range("A158").Value = "Islas Feroe"
$ws.Range("B158").Value = 187
$ws.Range("C158").Value = 0
$ws.Range("D158").Value = 187
$ws.Range("E158").Value = 0

# --- Row 160: Gibraltar ---
$ws.Range("B160").Value = 154
$ws.Range("C160").Value = 2
$ws.Range("E160").Value = 7
